$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.98"
$ws.Range("D3").Value = "'49.43"
$ws.Range("E3").Value = "'-0.52%"
$ws.Range("D4").Value = "'5.157"
$ws.Range("E4").Value = "'-3.37%"
$ws.Range("D5").Value = "'0.07732"
$ws.Range("E5").Value = "'-5.31%"
$ws.Range("D6").Value = "'4.520"
$ws.Range("E6").Value = "'-1.77%"
$ws.Range("D7").Value = "'1.366"
$ws.Range("E7").Value = "'12.86%"
$ws.Range("D8").Value = "'1.552"
$ws.Range("E8").Value = "'-7.25%"
$ws.Range("D9").Value = "'0.1229"
$ws.Range("E9").Value = "'-8.76%"
$ws.Range("D10").Value = "'0.1961"
$ws.Range("E10").Value = "'0.01%"
$ws.Range("D11").Value = "'0.09383"
$ws.Range("E11").Value = "'-2.84%"
$ws.Range("D12").Value = "'0.04657"
$ws.Range("E12").Value = "'5.56%"
$ws.Range("D13").Value = "'0.1045"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("E14").Value = "'-4.87%"
$ws.Range("D15").Value = "'0.04189"
$ws.Range("E15").Value = "'-2.63%"
$ws.Range("D16").Value = "'0.005841"
$ws.Range("E16").Value = "'0.42%"
$ws.Range("E17").Value = "'2,021.78%"
$ws.Range("E18").Value = "'-1.28%"
$ws.Range("D19").Value = "'2.239"
$ws.Range("E19").Value = "'-8.20%"
$ws.Range("D21").Value = "'7.995"
$ws.Range("E21").Value = "'-1.89%"
$ws.Range("D22").Value = "'0.1341"
$ws.Range("E22").Value = "'-5.51%"
$ws.Range("D23").Value = "'0.3041"
$ws.Range("E23").Value = "'4.43%"
$ws.Range("D24").Value = "'0.001275"
$ws.Range("E24").Value = "'-2.28%"
$ws.Range("D25").Value = "'0.004010"
$ws.Range("E25").Value = "'-6.07%"
$ws.Range("E26").Value = "'0.27%"
$ws.Range("D38").Value = "'0.02586"
$ws.Range("E38").Value = "'-6.35%"
$ws.Range("D39").Value = "'0.05818"
$ws.Range("E39").Value = "'3.94%"
$ws.Range("D40").Value = "'0.01075"
$ws.Range("E40").Value = "'70.71%"
$ws.Range("D41").Value = "'0.007909"
$ws.Range("E41").Value = "'2.85%"
$ws.Range("E42").Value = "'-2.00%"
$ws.Range("D43").Value = "'0.008465"
$ws.Range("E43").Value = "'10.27%"
$ws.Range("D44").Value = "'0.007703"
$ws.Range("E44").Value = "'-4.94%"
$ws.Range("D45").Value = "'0.3376"
$ws.Range("E45").Value = "'5.78%"
$ws.Range("D46").Value = "'0.00007030"
$ws.Range("E46").Value = "'1.02%"
$ws.Range("E47").Value = "'0.27%"
$ws.Range("D48").Value = "'0.05095"
$ws.Range("E48").Value = "'-16.92%"
$ws.Range("D49").Value = "'0.002627"
$ws.Range("E49").Value = "'-34.33%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.27%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.27%"
